# Generate Report for Handoff
#
# Moves the localization status from "In Translation" to "Ready for
# handoff" and stamps the new handoff-generation timestamps on the
# Overview sheet and on each per-language detail sheet (zh-cn, de-de).
# Excel's own column autosize then widens the "Status" column(s) to fit
# the new, longer text.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

# "Status" columns for zh-cn (E) and de-de (F)
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"

# "Latest HO Xliff Generate Date" (G)
$overview.Range("G2").Value = "2016-08-19 08:37:38"

# Widen the two Status columns to fit the new text (matches autofit).
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332

# ---- zh-cn detail sheet ------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "Ready for handoff"          # Status
$zhcn.Range("H2").Value = "2016-08-19 08:37:34"        # Latest Handoff Datetime
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332 # widen Status column

# ---- de-de detail sheet ------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Ready for handoff"          # Status
$dede.Range("H2").Value = "2016-08-19 08:37:38"        # Latest Handoff Datetime
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332 # widen Status column
